$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A43").Value = 45891
$ws.Range("A43").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B43").Value = 6466.91015625
$ws.Range("C43").Value = 6478.89013671875
$ws.Range("D43").Value = 6384.58984375
$ws.Range("E43").Value = 6384.58984375
$ws.Range("F43").Value = 5046860000
$ws.Range("G43").Value = 0.0151864448768934
$ws.Range("H43").Value = 1

$ws.Range("I43").Value = 'blockchain biweeklysec issue guidance regarding liquid staking token and trial of privacy technology developer end in mixed verdict forbes name top employer in arkansas for nasdaq double top near aligns with gann timing for cycle peak dr phil accused of fraud in m tv deal gone wrong coalinga unit property sell for m setting city record the doj is considering a probe into fed lisa cook over mortgage fraud allegation home depot labor day sale start now here everything im buying for fall google add agentic feature to aipowered search essential financial tool every small business should know in essential financial tool every small business should know in walmart post revenue growth in q raise fiscal sale guidance finding their voice ccisd look to alternative amid speechlanguage pathologist shortage sweet southern sas in seabrook reopens after renovation shiba inu shib price climb a analyst predict it could erase another zero while one new meme coin target a run piada italian street food coming to pearland town center ny appeal court throw out million civil fraud penalty against trump the wall street journal ukrainian arrested in italy for nord stream pipeline explosion dogecoin bonk and magax which meme coin could dominate the next bull run vantage drilling international ltd further extension of conditional letter of award vantage drilling international ltd further extension of conditional letter of award california is too overregulated expensive and risky for bed bath beyond chairman say why coty stock is plummeting today air canada employee accelerate the recovery airline update goodwill policy for customer to exceptionally include accommodation and outofpocket cost walmart stock fall on mixed q result despite raised guidance ionis pharma win fda nod for first rnatargeted therapy for rare genetic swelling disorder extra extra there are just so many interesting way to acquire money accredited wealth management llc lower position in exxon mobil corporation xom housing development with a lazy river fresno county town greenlights resortstyle living the surprising reason leaving the office early is fine but showing up late is not willis isd approves m balanced budget for fy bilt alaska atmos reward partnership enticing new opportunity one mile at a time silver rally aligns with gann geometry but face august crest risk the gdl fund declares third quarter distribution of per share will county executive jennifer bertinotarrant boast of success present million budget edenville dam restoration take shape a largest four lake undertaking investor alert pomerantz law firm investigates claim on behalf of investor of agilon health the gdl fund declares third quarter distribution of per share fund finance review bad credit personal loan in a layoff and inflation push household toward emergency funding new york appeal court toss m civil fraud penalty against trump al jazeera make bedside screen time more convenient with the magicoutlet nc trial court reject sealing of part of state complaint against tiktok porter county homeowner ready to catch some ray under new personal solar ordinance free spin no deposit casino uptown ace launch free spin bonus new fox one streamer launching on prime video deadline top new meme coin to invest in now that investor are calling the next big bet i just found these smartphones with inch display and it amazing justice department urge fed chair to fire top official accused of mortgage fraud alien craft or space rock wall street race to lift nvidia price target ahead of earnings bloombergcom u flight attendant are fed up like their air canada peer here why they are unlikely to strike achieve life science inch closer to fda nod for first new quitsmoking drug in year tecnoglass responds to short seller report u manufacturing is picking up but so are consumer price county supes want to tap reerve to offset effect of one big beautiful bill walmarts revenue rise profit growth fall short of expectation palantir is selling off bigbearai could be your hedge new facility and interest in watertown international airport the link between hypercompetition and freight fraud china firm launch solar panel steering battery that work even in extreme f ny appeal court dump trump million fraud fine but upholds judgment gabelli global utility income trust continues monthly distribution declares distribution of per share gabelli global utility income trust continues monthly distribution declares distribution of per share u manufacturing surge to year high a demand strengthens hiring pick up useu trade deal new framework agreement boost transatlantic economic tie castlegar award contract for based on previous track record castlegar are mouse fruit fly the new space crew ethereum eth surpasses bitcoins btc yeartodate gain while market shift attention to viral new crypto for roi sharkninja option trading a deep dive into market sentiment deep dive into charter communication stock analyst perspective rating unpacking the latest option trading trend in shake shack this is what whale are betting on pfizer assessing hilton grand vacation insight from financial analyst if bitcoin treasury company are buying billion why is the price not going up breaking down paychex analyst share their view cracker barrel share plummet after pushback on new logo brand refresh cnbc demystifying insulet insight from analyst review deadline alert cto realty growth inc cto investor who lost money urged to contact glancy prongay murray llp about security fraud lawsuit subnautica s developer unknown world sue founder a krafton lawsuit saga continues mortenson open advanced research facility in fridley the analyst verdict sociedad quimica y minera in the eye of expert pe ratio insight for magnite is progressive gaining or losing market support i bought worth of cheap home decor from amazon haul here what id buy again looking into johnson control intls recent short interest dexcom stock a deep dive into analyst perspective rating how do investor really feel about johnson johnson whyy newscast for thursday pm orange county football week storyline a mix of oldest and youngest deadline alert centene corporation cnc investor who lost money urged to contact glancy prongay murray llp about security fraud lawsuit earnings season reveals winner and straggler across retail offprice arena constableville news west turin review plan for new sandsalt building xrp news today xrp drop below a sec postpones etf ruling until october teleflex launch barrigeltm rectal spacer in japan following regulatory approval northeast community bancorp inc announces special cash dividend oklahoma voice oklahoma tset board hope to pressure some major company to change their policy northeast community bancorp inc announces special cash dividend tegna name brent denny president and general manager at wthr in indianapolis near eye cro forecast gain but cold wallet selfcustody set new decentralized standard silo pharma alzheimers candidate spc featured in leading scientific peerreviewed journal u national security probe target wind industry orlando official denounce removal of rainbow crosswalk near pulse nightclub mass shooting site orlando official denounce removal of rainbow crosswalk near pulse nightclub mass shooting site orlando official denounce removal of rainbow crosswalk near pulse nightclub mass shooting site orlando official denounce removal of rainbow crosswalk near pulse nightclub mass shooting site swingstate gov slam trump plan to nix mailin voting he can sign whatever the hell he want trumpdecreed smithsonian scrubbing humbles oncegreat nation former speaker dade phelan wont seek reelection to texas house trump administration use of travel ban dealt blow from judge elderly man tricked into signing over deed to home investigator said trump wont stop deporting immigrant to country they arent from commentary defending the enlightenment ii south korea plan to stop north korea nuclear weapon program city of norman emphasize work zone safety after recent close call what trump crime crackdown in dc look like in photo senate democrat predicts day of reckoning for private prison employee tinley park district teacher union declares impasse over contract negotiation fatal florida crash fuel immigration fight between desantis and newsom trump dumbest trade adviser attack fed chair for dumbest reason carney ha productive phone call with trump amid bilateral tension canada florida turnpike uturn crash that killed fuel immigration fight between govs desantis and newsom lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention lawyer ask judge to order ice to free spanishlanguage journalist from immigration detention trump plan to patrol street of u capital with troop trump plan to patrol street of u capital with troop trump plan to patrol street of u capital with troop trump plan to patrol street of u capital with troop uscis chooses for agile development pact ap news summary at pm edt antiice protester speaks out after arrest in san francisco demoted general sue utah national guard governor over firing netanyahu say ready for negotiation a gaza city takeover advance effective immediately state department halt worker visa for commercial truck driver a week of summit reveals trump is closer to moscow than europe city solicitor did not force cincinnati police to charge white man in violent fight despite fop claim u issue sanction against greek shipper of iranian oil two baby hospitalized sewage flooding home family sick lt governor beckwith and alexandria official silent say concerned citizen of alexandria is the redistricting fight abortion andy mccarthy explains why appeal court tossed letitia james ridiculous penalty against trump huge win for donald trump cnns elie honig explains monumental ruling against letitia james federal funding for sex education in california is cut over radical gender ideology orphan child allegedly bribed into voting for btss jimin sparking viral outrage u citizen legal resident detained during socal immigration raid filing claim against government israel pm order talk to free all our hostage a army pound gaza city israel pm order talk to free all our hostage a army pound gaza city israel pm order talk to free all our hostage a army pound gaza city israel pm order talk to free all our hostage a army pound gaza city israel pm order talk to free all our hostage a army pound gaza city most say partisan redistricting threatens democracy survey house oversight chair say panel to release epstein file a quickly a we can while protecting victim which country ha the most nuclear bunker new window driving law spark wild west concern cop can still fine you even after check removed fact focus a look at rfk jr misleading claim on u dietary guideline and froot loop fact focus a look at rfk jr misleading claim on u dietary guideline and froot loop fact focus a look at rfk jr misleading claim on u dietary guideline and froot loop fact focus a look at rfk jr misleading claim on u dietary guideline and froot loop fact focus a look at rfk jr misleading claim on u dietary guideline and froot loop fact focus a look at rfk jr misleading claim on u dietary guideline and froot loop fact focus a look at rfk jr misleading claim on u dietary guideline and froot loop fact focus a look at rfk jr misleading claim on u dietary guideline and froot loop fact focus a look at rfk jr misleading claim on u dietary guideline and froot loop california city lack unified response on homeless encampment exgop leader go on cnn immediately continues feud with trump loyalist a nonprofit want to resume the search for the missing cabazon baby and is baffled the parent have not responded trump demand release of colorado official sentenced in election data scheme the washington post bomb threat at north kamloops ministry office investigated by bcgeu occupational health and safety officer jasmine crockett disparages maga voter trump a most unpatriotic people hhs terminates california sex education grant over gender ideology content ucla in the news august colorado lawmaker seek censure of republican who shared photo of democrat prompting mocking comment colorado lawmaker seek censure of republican who shared photo of democrat prompting mocking comment editorial after a disastrous alaska summit a feeble trump look poised to complete his sellout of ukraine editorial after a disastrous alaska summit a feeble trump look poised to complete his sellout of ukraine orlando official denounce removal of rainbow crosswalk near pulse nightclub mass shooting site trump official called protester against national guard in dc stupid white hippy state lead u with most unauthorized immigrant new pew study find axios belly acre burger joint evicted from east memphis regalia shopping center what to know netanyahu say hell push ahead with gaza city takeover and renewed ceasefire talk netanyahu say hell push ahead with gaza city takeover and renewed ceasefire talk netanyahu say hell push ahead with gaza city takeover and renewed ceasefire talk netanyahu say hell push ahead with gaza city takeover and renewed ceasefire talk ap news summary at pm edt netanyahu say hell push ahead with gaza city takeover and renewed ceasefire talk'
$ws.Range("J43").Value = 'blockchain biweeklysec issue guidance regarding liquid staking token trial privacy technology developer end mixed verdict forbes name employer arkansas nasdaq double aligns gann timing cycle peak dr phil accused fraud tv deal coalinga unit property sell setting record doj considering probe fed lisa cook mortgage fraud allegation depot labor day sale start everything im buying fall google add agentic feature aipowered search essential financial tool small business essential financial tool small business walmart post revenue growth q raise fiscal sale guidance finding voice ccisd look alternative amid speechlanguage pathologist shortage sweet southern sa seabrook reopens renovation shiba inu shib price climb analyst predict erase another new meme coin target run piada italian street food coming pearland town center ny appeal court throw million civil fraud penalty trump wall street journal ukrainian arrested italy nord stream pipeline explosion dogecoin bonk magax meme coin dominate next bull run vantage drilling international extension conditional letter award vantage drilling international extension conditional letter award california overregulated expensive risky bed bath beyond chairman coty stock plummeting today air canada employee accelerate recovery airline update goodwill policy customer exceptionally include accommodation outofpocket cost walmart stock fall mixed q raised guidance ionis pharma fda nod rnatargeted therapy rare genetic swelling disorder extra extra interesting acquire money accredited wealth management llc lower exxon mobil corporation xom housing development lazy river fresno county town greenlights resortstyle living surprising leaving office early fine showing late willis isd approves balanced budget fy bilt alaska atmos reward partnership enticing new opportunity mile silver rally aligns gann geometry face august crest risk gdl fund declares third quarter distribution per share county executive jennifer bertinotarrant boast success present million budget edenville dam restoration take shape largest lake undertaking investor alert pomerantz law firm investigates claim behalf investor agilon health gdl fund declares third quarter distribution per share fund finance review bad credit personal loan layoff inflation push household emergency funding new york appeal court toss civil fraud penalty trump al jazeera bedside screen convenient magicoutlet nc trial court reject sealing part complaint tiktok porter county homeowner ready catch ray new personal solar ordinance free spin deposit casino uptown ace launch free spin bonus new fox streamer launching prime video deadline new meme coin invest investor calling next big bet found smartphones inch display amazing justice department urge fed chair fire accused mortgage fraud alien craft space rock wall street race lift nvidia price target ahead earnings bloombergcom flight attendant fed air canada peer unlikely strike achieve life science inch closer fda nod new quitsmoking drug year tecnoglass responds short seller report manufacturing picking consumer price county supes tap reerve offset effect big bill walmarts revenue rise profit growth fall short expectation palantir selling bigbearai hedge new facility interest watertown international airport link hypercompetition freight fraud china firm launch solar panel steering battery work extreme f ny appeal court dump trump million fraud fine upholds judgment gabelli global utility income trust continues monthly distribution declares distribution per share gabelli global utility income trust continues monthly distribution declares distribution per share manufacturing surge year high demand strengthens hiring pick useu trade deal new framework agreement boost transatlantic economic tie castlegar award contract previous track record castlegar mouse fruit fly new space crew ethereum eth surpasses bitcoins btc yeartodate gain market shift attention viral new crypto roi sharkninja option trading deep dive market sentiment deep dive charter communication stock analyst perspective rating unpacking latest option trading trend shake shack whale betting pfizer assessing hilton grand vacation insight financial analyst bitcoin treasury company buying billion price breaking paychex analyst share view cracker barrel share plummet pushback new logo brand refresh cnbc demystifying insulet insight analyst review deadline alert cto realty growth cto investor money urged contact glancy prongay murray llp security fraud lawsuit subnautica developer unknown world sue founder krafton lawsuit saga continues mortenson open advanced research facility fridley analyst verdict sociedad quimica minera eye expert pe ratio insight magnite progressive gaining losing market support bought worth cheap decor amazon haul id buy johnson control intls recent short interest dexcom stock deep dive analyst perspective rating investor feel johnson johnson whyy newscast thursday pm orange county football storyline mix oldest youngest deadline alert centene corporation cnc investor money urged contact glancy prongay murray llp security fraud lawsuit earnings season reveals winner straggler across retail offprice arena constableville news west turin review plan new sandsalt building xrp news today xrp drop sec postpones etf ruling october teleflex launch barrigeltm rectal spacer japan following regulatory approval northeast community bancorp announces special cash dividend oklahoma voice oklahoma tset board pressure major company change policy northeast community bancorp announces special cash dividend tegna name brent denny president general manager wthr indianapolis eye cro forecast gain cold wallet selfcustody set new decentralized silo pharma alzheimers candidate spc featured leading scientific peerreviewed journal national security probe target wind industry orlando denounce removal rainbow crosswalk pulse nightclub mass shooting site orlando denounce removal rainbow crosswalk pulse nightclub mass shooting site orlando denounce removal rainbow crosswalk pulse nightclub mass shooting site orlando denounce removal rainbow crosswalk pulse nightclub mass shooting site swingstate gov slam trump plan nix mailin voting sign whatever hell trumpdecreed smithsonian scrubbing humbles oncegreat nation speaker dade phelan wont seek reelection texas house trump administration travel ban dealt blow judge elderly man tricked signing deed investigator trump wont stop deporting immigrant country arent commentary defending enlightenment ii south korea plan stop north korea nuclear weapon program norman emphasize work zone safety recent close call trump crime crackdown dc look photo senate democrat predicts day reckoning private prison employee tinley park district teacher union declares impasse contract negotiation fatal florida crash fuel immigration fight desantis newsom trump dumbest trade adviser attack fed chair dumbest carney productive phone call trump amid bilateral tension canada florida turnpike uturn crash killed fuel immigration fight govs desantis newsom lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention lawyer judge order free spanishlanguage journalist immigration detention trump plan patrol street capital troop trump plan patrol street capital troop trump plan patrol street capital troop trump plan patrol street capital troop uscis chooses agile development pact news summary pm antiice protester speaks arrest san francisco demoted general sue utah national guard governor firing netanyahu ready negotiation gaza takeover advance effective immediately department halt worker visa commercial truck driver summit reveals trump closer moscow europe solicitor force cincinnati police charge white man violent fight fop claim issue sanction greek shipper iranian oil baby hospitalized sewage flooding family sick lt governor beckwith alexandria silent concerned citizen alexandria redistricting fight abortion andy mccarthy explains appeal court tossed letitia james ridiculous penalty trump donald trump cnns elie honig explains monumental ruling letitia james federal funding sex education california cut radical gender ideology orphan child allegedly bribed voting btss jimin sparking viral outrage citizen legal resident detained socal immigration raid filing claim government israel pm order free hostage army pound gaza israel pm order free hostage army pound gaza israel pm order free hostage army pound gaza israel pm order free hostage army pound gaza israel pm order free hostage army pound gaza partisan redistricting threatens democracy survey house oversight chair panel release epstein file quickly protecting victim country nuclear bunker new window driving law spark wild west concern cop fine removed fact focus look rfk jr misleading claim dietary guideline froot loop fact focus look rfk jr misleading claim dietary guideline froot loop fact focus look rfk jr misleading claim dietary guideline froot loop fact focus look rfk jr misleading claim dietary guideline froot loop fact focus look rfk jr misleading claim dietary guideline froot loop fact focus look rfk jr misleading claim dietary guideline froot loop fact focus look rfk jr misleading claim dietary guideline froot loop fact focus look rfk jr misleading claim dietary guideline froot loop fact focus look rfk jr misleading claim dietary guideline froot loop california lack unified response homeless encampment exgop cnn immediately continues feud trump loyalist nonprofit resume search missing cabazon baby baffled parent responded trump demand release colorado sentenced election data scheme washington post bomb threat north kamloops ministry office investigated bcgeu occupational health safety officer jasmine crockett disparages maga voter trump unpatriotic people hhs terminates california sex education grant gender ideology content ucla news august colorado lawmaker seek censure republican shared photo democrat prompting mocking colorado lawmaker seek censure republican shared photo democrat prompting mocking editorial disastrous alaska summit feeble trump look poised complete sellout ukraine editorial disastrous alaska summit feeble trump look poised complete sellout ukraine orlando denounce removal rainbow crosswalk pulse nightclub mass shooting site trump protester national guard dc stupid white hippy lead unauthorized immigrant new pew study find axios belly acre burger joint evicted east memphis regalia shopping center netanyahu hell push ahead gaza takeover renewed ceasefire netanyahu hell push ahead gaza takeover renewed ceasefire netanyahu hell push ahead gaza takeover renewed ceasefire netanyahu hell push ahead gaza takeover renewed ceasefire news summary pm netanyahu hell push ahead gaza takeover renewed ceasefire'
